# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 6, pushing the existing
# rows 6-14 down to rows 7-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 6..14 down by one row (Excel-style row insert, like
# right-clicking row 6 and choosing "Insert").
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new weekly record.
$ws.Range("A6").Value = 12
$ws.Range("B6").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C6").Value = "Metropolitana"
$ws.Range("D6").Value = 45030
$ws.Range("D6").NumberFormat = $ws.Range("D7").NumberFormat
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = 100112032
$ws.Range("G6").Value = "Zapallo italiano"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 6000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 6000
$ws.Range("N6").Value = "`$/caja 50 unidades"
$ws.Range("O6").Value = "Región de Arica y Parinacota"
$ws.Range("P6").Value = 120
$ws.Range("Q6").Value = 50
$ws.Range("R6").Value = "Hortaliza"
